{"js": "// The diff does two things inside this r\u00e9sum\u00e9 document:\n//  1. In the paragraph \"2023\u5e74<TAB>3\u670823\u65e5\", the run that held both the\n//     tab character and the text \"3\u670823\u65e5\" is split into two runs (same\n//     run properties on both): one run keeping just the <w:tab/>, and a\n//     brand-new run holding the text \"3\u670823\u65e5\".\n//  2. In the paragraph whose text is exactly \"\u4eca\u5929\u661f\u671f\u56db\", \" 123\" is\n//     appended to the end of the text (so it becomes \"\u4eca\u5929\u661f\u671f\u56db 123\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n\n  // --- Change 1: split the \"<tab>3\u670823\u65e5\" run into two runs ---\n  if (text.indexOf(\"2023\") !== -1 && text.indexOf(\"3\\u670823\\u65e5\") !== -1 && text.indexOf(\"\\t\") !== -1) {\n    const range = paragraph.getRange();\n    const ooxml = paragraph.getOoxml();\n    await context.sync();\n\n    // Pull the paragraph-mark run properties (<w:pPr><w:rPr>) and the run\n    // properties used by the run we are about to split, straight out of the\n    // paragraph's own OOXML, so the split keeps formatting byte-identical.\n    const raw = ooxml.value;\n    const pPrMatch = raw.match(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/);\n    const pPrXml = pPrMatch ? pPrMatch[0] : \"\";\n    const runRPrMatch = raw.match(/<w:r><w:rPr>[\\s\\S]*?<\\/w:rPr><w:tab\\/>/);\n    const rPrXml = runRPrMatch ? runRPrMatch[0].match(/<w:rPr>[\\s\\S]*?<\\/w:rPr>/)[0] : \"<w:rPr/>\";\n    const firstRunMatch = raw.match(/<w:r>(?:(?!<w:r>)[\\s\\S])*?<w:t>2023[\\s\\S]*?<\\/w:t><\\/w:r>/);\n    const firstRunXml = firstRunMatch ? firstRunMatch[0] : `<w:r>${rPrXml}<w:t>2023\\u5e74</w:t></w:r>`;\n\n    const newParagraphOoxml =\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body>' +\n      '<w:p>' + pPrXml +\n      firstRunXml +\n      '<w:r>' + rPrXml + '<w:tab/></w:r>' +\n      '<w:r>' + rPrXml + '<w:t>3\\u670823\\u65e5</w:t></w:r>' +\n      '</w:p>' +\n      '</w:body></w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>';\n\n    range.insertOoxml(newParagraphOoxml, \"Replace\");\n    await context.sync();\n  }\n\n  // --- Change 2: append \" 123\" to the \"\u4eca\u5929\u661f\u671f\u56db\" paragraph ---\n  if (text === \"\\u4eca\\u5929\\u661f\\u671f\\u56db\") {\n    paragraph.insertText(\" 123\", \"End\");\n    await context.sync();\n  }\n}\n", "ps1": "# The diff does two things inside this r\u00e9sum\u00e9 document:\n#  1. In the paragraph \"2023\u5e74<TAB>3\u670823\u65e5\", the run that held both the\n#     tab character and the text \"3\u670823\u65e5\" is split into two runs (same\n#     run properties on both): one run keeping just the <w:tab/>, and a\n#     brand-new run holding the text \"3\u670823\u65e5\".\n#  2. In the paragraph whose text is exactly \"\u4eca\u5929\u661f\u671f\u56db\", \" 123\" is\n#     appended to the end of the text (so it becomes \"\u4eca\u5929\u661f\u671f\u56db 123\").\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text\n    $trimmed = $paraText.TrimEnd([char]13, [char]7)\n\n    # --- Change 1: split the \"<tab>3\u670823\u65e5\" run into two runs ---\n    if ($trimmed.Contains(\"2023\") -and $trimmed.Contains([char]9)) {\n        $r = $p.Range\n        $oldXml = $r.WordOpenXML\n\n        # Everything from the start of this <w:p> up to (not including) the\n        # run that carries the <w:tab/> - i.e. the <w:pPr> plus the \"2023\u5e74\"\n        # run, copied verbatim so formatting stays identical.\n        $bodyStart = $oldXml.IndexOf(\"<w:body>\") + 8\n        $pStart = $oldXml.IndexOf(\"<w:p\", $bodyStart)\n        $pContentStart = $oldXml.IndexOf(\">\", $pStart) + 1\n\n        $tabIdx = $oldXml.IndexOf(\"<w:tab/>\")\n        $tabRunStart = $oldXml.LastIndexOf(\"<w:r>\", $tabIdx)\n\n        $beforeXml = $oldXml.Substring($pContentStart, $tabRunStart - $pContentStart)\n\n        # Run properties of the run that contains the tab (reused verbatim\n        # for both halves of the split).\n        $rPrOpenTag = $oldXml.IndexOf(\"<w:rPr>\", $tabRunStart)\n        $rPrCloseTag = $oldXml.IndexOf(\"</w:rPr>\", $rPrOpenTag) + 8\n        $rPrXml = $oldXml.Substring($rPrOpenTag, $rPrCloseTag - $rPrOpenTag)\n\n        # Whatever followed the <w:tab/> inside that run (the \"3\u670823\u65e5\" text).\n        $afterTabIdx = $tabIdx + 8\n        $runEndIdx = $oldXml.IndexOf(\"</w:r>\", $afterTabIdx)\n        $afterTabXml = $oldXml.Substring($afterTabIdx, $runEndIdx - $afterTabIdx)\n\n        $newXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $beforeXml + '<w:r>' + $rPrXml + '<w:tab/></w:r><w:r>' + $rPrXml + $afterTabXml + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n        $r.InsertXML($newXml)\n    }\n\n    # --- Change 2: append \" 123\" to the \"\u4eca\u5929\u661f\u671f\u56db\" paragraph ---\n    if ($trimmed -eq \"\u4eca\u5929\u661f\u671f\u56db\") {\n        $r2 = $p.Range\n        $r2.MoveEnd(1, -1) | Out-Null\n        $r2.InsertAfter(\" 123\")\n    }\n}\n"}
